# Update the Coin / Link / Price / Volume(1h) figures for the cryptos list
# to match the latest refresh pulled from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a literal text value to a cell without letting Excel
# reinterpret numeric/percentage-looking strings as numbers, and without
# leaving the cell with a different number format / style than before.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 '44.011.01'
Set-TextValue 2 5 '  +2.33%  '

# Row 3
Set-TextValue 3 4 '2.256.42'
Set-TextValue 3 5 '  +1.59%  '

# Row 4
Set-TextValue 4 5 '  +0.02%  '

# Row 5
Set-TextValue 5 4 '319.55'
Set-TextValue 5 5 '  -0.38%  '

# Row 6
Set-TextValue 6 4 '102.13'
Set-TextValue 6 5 '  +3.15%  '

# Row 7
Set-TextValue 7 5 '  -0.59%  '

# Row 8
Set-TextValue 8 5 '  +0.12%  '

# Row 9
Set-TextValue 9 4 '0.555'
Set-TextValue 9 5 '  -0.30%  '

# Row 10
Set-TextValue 10 4 '37.52'
Set-TextValue 10 5 '  +2.09%  '

# Row 11
Set-TextValue 11 4 '0.0833'
Set-TextValue 11 5 '  +0.83%  '

# Row 12
Set-TextValue 12 4 '7.65'
Set-TextValue 12 5 '  +0.36%  '

# Row 13
Set-TextValue 13 5 '  -1.36%  '

# Row 14
Set-TextValue 14 4 '2.603.05'
Set-TextValue 14 5 '  +1.62%  '

# Row 15
Set-TextValue 15 4 '0.861'
Set-TextValue 15 5 '  -0.18%  '

# Row 16
Set-TextValue 16 4 '14.50'
Set-TextValue 16 5 '  +1.07%  '

# Row 17
Set-TextValue 17 4 '2.261.45'
Set-TextValue 17 5 '  +1.85%  '

# Row 18
Set-TextValue 18 4 '43.956.11'
Set-TextValue 18 5 '  +2.34%  '

# Row 19
Set-TextValue 19 4 '13.47'
Set-TextValue 19 5 '  -4.22%  '

# Row 20
Set-TextValue 20 4 '0.0₃0988'
Set-TextValue 20 5 '  +2.75%  '

# Row 21
Set-TextValue 21 4 '6.54'
Set-TextValue 21 5 '  -0.01%  '

# Row 22
Set-TextValue 22 4 '65.90'
Set-TextValue 22 5 '  +1.48%  '

# Row 23
Set-TextValue 23 4 '3.16'

# Row 24
Set-TextValue 24 4 '235.68'
Set-TextValue 24 5 '  -0.22%  '

# Row 25
Set-TextValue 25 5 '  -2.13%  '

# Row 26
Set-TextValue 26 5 '  +0.10%  '

# Row 27
Set-TextValue 27 5 '  +2.95%  '

# Row 28
Set-TextValue 28 4 '38.00'
Set-TextValue 28 5 '  +6.79%  '

# Row 29
Set-TextValue 29 4 '2.20'
Set-TextValue 29 5 '  -2.03%  '

# Row 30
Set-TextValue 30 4 '6.24'
Set-TextValue 30 5 '  -1.29%  '

# Row 31
Set-TextValue 31 4 '161.16'
Set-TextValue 31 5 '  +5.87%  '

# Row 32
Set-TextValue 32 4 '20.23'
Set-TextValue 32 5 '  +0.01%  '

# Row 33
Set-TextValue 33 4 '0.0853'
Set-TextValue 33 5 '  -1.50%  '

# Row 34
Set-TextValue 34 5 '  +0.34%  '

# Row 35
Set-TextValue 35 5 '  +11.19%  '

# Row 36
Set-TextValue 36 5 '  +2.49%  '

# Row 37
Set-TextValue 37 4 '3.05'
Set-TextValue 37 5 '  -2.51%  '

# Row 38
Set-TextValue 38 5 '  -1.87%  '

# Row 39
Set-TextValue 39 4 '16.71'
Set-TextValue 39 5 '  +22.00%  '

# Row 40
Set-TextValue 40 4 '3.74'
Set-TextValue 40 5 '  +2.25%  '

# Row 41
Set-TextValue 41 5 '  -4.09%  '

# Row 42
Set-TextValue 42 5 '  -1.28%  '

# Row 43
Set-TextValue 43 5 '  +0.22%  '

# Row 44
Set-TextValue 44 4 '1.801.35'

# Row 45 (Coin -> ordi)
Set-TextValue 45 2 'ordi'
Set-TextValue 45 3 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue 45 4 '76.34'
Set-TextValue 45 5 '  +2.82%  '

# Row 46 (Coin -> Algorand)
Set-TextValue 46 2 'Algorand'
Set-TextValue 46 3 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 46 4 '0.199'
Set-TextValue 46 5 '  -1.89%  '

# Row 47
Set-TextValue 47 4 '82.83'
Set-TextValue 47 5 '  -1.98%  '

# Row 48
Set-TextValue 48 4 '5.23'
Set-TextValue 48 5 '  -0.78%  '

# Row 49
Set-TextValue 49 4 '104.95'
Set-TextValue 49 5 '  +1.81%  '

# Row 50
Set-TextValue 50 5 '  +8.42%  '

# Row 51
Set-TextValue 51 4 '58.42'
Set-TextValue 51 5 '  +1.04%  '
